# Apply "nuevos experimentos no convexos" updates to the generator workbook.
# Updates numeric/text experiment data across several sheets. Many of the
# cells hold *text* representations of numbers (they were produced by a
# generator script, not typed into Excel), so for pure-numeric-looking
# strings we briefly force Text number format before assigning the value
# (otherwise Excel auto-converts the string to a real number), then clear
# the format again so the cell itself carries no leftover styling.

$wb = $excel.ActiveWorkbook

# --- Sheet "Restricciones_del_follower" (3rd sheet) ---------------------
$ws3 = $wb.Worksheets.Item(3)

# Force any numeric-looking text in this block to stay text.
$ws3.Range("B2:B4").NumberFormat = "@"
$ws3.Range("D2:F4").NumberFormat = "@"

# Row 2 (J_0_L0_v)
$ws3.Range("A2").Value = "-3.5801944728761512 - 2x_1 + 1.5730467417263732y_1 + 1.9255544182872733y_2"
$ws3.Range("B2").Value = "6.080194472876151"
$ws3.Range("F2").Value = "6.5"

# Row 3 (J_0_LP_v)
$ws3.Range("A3").Value = "1.3934237461617203 + x_1 - 3x_2 - 0.18863868986693957y_1 + 0.03694984646878219y_2"
$ws3.Range("B3").Value = "-3.3934237461617203"
$ws3.Range("D3").Value = "0.7"
$ws3.Range("E3").Value = "0"
$ws3.Range("F3").Value = "9.9"

# Row 4 (J_Ne_L0_v)
$ws3.Range("A4").Value = "-6.6 + x_1 + x_2"
$ws3.Range("B4").Value = "4.1"
$ws3.Range("D4").Value = "0.8"
$ws3.Range("E4").Value = "2.9"
$ws3.Range("F4").Value = "0"

$ws3.Range("B2:B4").ClearFormats()
$ws3.Range("D2:F4").ClearFormats()

# --- Sheet "Punto_modificado" (4th sheet) --------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2:D2").NumberFormat = "@"
$ws4.Range("A2").Value = "4.5"
$ws4.Range("B2").Value = "1.6"
$ws4.Range("C2").Value = "6.1000000000000005"
$ws4.Range("D2").Value = "1.55"
$ws4.Range("A2:D2").ClearFormats()

# --- Sheet "Vector_bf" (5th sheet) ---------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2:A3").NumberFormat = "@"
$ws5.Range("A2").Value = "2.6219222108495392"
$ws5.Range("A3").Value = "-2.87439713408393"
$ws5.Range("A2:A3").ClearFormats()

# --- Sheet "Vector_BF" (6th sheet) ---------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2:A5").NumberFormat = "@"
$ws6.Range("A2").Value = "-0.8999999999999999"
$ws6.Range("A3").Value = "-3.9"
$ws6.Range("A4").Value = "-0.5"
$ws6.Range("A5").Value = "-0.0"
$ws6.Range("A2:A5").ClearFormats()

# --- Sheet "Vector_Alpha" (7th sheet) ------------------------------------
# These two stay real numeric cells (no t="s" in the source file).
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 0.5700000000000001
$ws7.Range("A3").Value = 2.91
